$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.922.68'
$ws.Range('E2').Value = '  +0.80%  '

$ws.Range('D3').Value = '1.636.21'
$ws.Range('E3').Value = '  +1.96%  '

$ws.Range('E4').Value = '  +0.27%  '

$ws.Range('D5').Value = '215.08'
$ws.Range('E5').Value = '  +1.10%  '

$ws.Range('D6').Value = '0.517'
$ws.Range('E6').Value = '  +0.06%  '

$ws.Range('E7').Value = '  +0.35%  '

$ws.Range('D8').Value = '28.98'
$ws.Range('E8').Value = '  +3.47%  '

$ws.Range('E9').Value = '  +2.34%  '

$ws.Range('E10').Value = '  +0.74%  '

$ws.Range('D11').Value = '0.0916'

$ws.Range('D12').Value = '1.872.43'
$ws.Range('E12').Value = '  +2.14%  '

$ws.Range('D13').Value = '1.634.37'
$ws.Range('E13').Value = '  +1.89%  '

$ws.Range('D14').Value = '0.567'
$ws.Range('E14').Value = '  +2.92%  '

$ws.Range('D15').Value = '9.34'
$ws.Range('E15').Value = '  +19.32%  '

$ws.Range('D16').Value = '3.89'
$ws.Range('E16').Value = '  +3.35%  '

$ws.Range('D17').Value = '29.906.84'
$ws.Range('E17').Value = '  +0.77%  '

$ws.Range('D18').Value = '64.32'
$ws.Range('E18').Value = '  +0.44%  '

$ws.Range('D19').Value = '243.93'
$ws.Range('E19').Value = '  +0.52%  '

$ws.Range('D20').Value = '0.0₃0703'
$ws.Range('E20').Value = '  +0.69%  '

$ws.Range('E21').Value = '  +0.28%  '

$ws.Range('D22').Value = '9.91'
$ws.Range('E22').Value = '  +5.23%  '

$ws.Range('E23').Value = '  +2.92%  '

$ws.Range('E24').Value = '  +1.70%  '

$ws.Range('D25').Value = '158.15'
$ws.Range('E25').Value = '  +1.74%  '

$ws.Range('D26').Value = '15.59'
$ws.Range('E26').Value = '  +0.68%  '

$ws.Range('D27').Value = '0.111'
$ws.Range('E27').Value = '  +1.77%  '

$ws.Range('E28').Value = '  +2.80%  '

$ws.Range('E29').Value = '  +0.47%  '

$ws.Range('D30').Value = '0.0489'
$ws.Range('E30').Value = '  +1.23%  '

$ws.Range('E31').Value = '  +5.24%  '

$ws.Range('E32').Value = '  +4.47%  '

$ws.Range('D33').Value = '3.18'
$ws.Range('E33').Value = '  -0.33%  '

$ws.Range('D34').Value = '1.429.50'
$ws.Range('E34').Value = '  -0.09%  '

$ws.Range('E35').Value = '  +6.22%  '

$ws.Range('E36').Value = '  +1.00%  '

$ws.Range('D37').Value = '2.84'
$ws.Range('E37').Value = '  -2.95%  '

$ws.Range('E38').Value = '  +0.50%  '

$ws.Range('B39').Value = 'Aave'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D39').Value = '77.34'
$ws.Range('E39').Value = '  +16.36%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0170'
$ws.Range('E40').Value = '  +0.78%  '

$ws.Range('D41').Value = '0.554'
$ws.Range('E41').Value = '  +1.09%  '

$ws.Range('D42').Value = '2.01'
$ws.Range('E42').Value = '  +2.00%  '

$ws.Range('D43').Value = '0.831'
$ws.Range('E43').Value = '  +1.59%  '

$ws.Range('E44').Value = '  -0.96%  '

$ws.Range('E45').Value = '  +5.82%  '

$ws.Range('D46').Value = '53.88'
$ws.Range('E46').Value = '  -7.07%  '

$ws.Range('E47').Value = '  +0.43%  '

$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '1.778.73'
$ws.Range('E48').Value = '  +2.18%  '

$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '5.36'
$ws.Range('E49').Value = '  +0.19%  '

$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '89.63'
$ws.Range('E50').Value = '  +3.28%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0112'
$ws.Range('E51').Value = '  +6.67%  '
